$wb = $excel.ActiveWorkbook

# Sheet "Overview": Latest HO Xliff Generate Date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-27 11:04:16"

# Sheet "zh-cn": Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-27 11:04:11"
$wsZhCn.Range("K2").Value = "2016-08-27 11:04:33"

# Sheet "de-de": Latest HO Xliff Generate Date / Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-27 11:04:16"
$wsDeDe.Range("K2").Value = "2016-08-27 11:04:39"
